$p = $ppt.ActivePresentation

# --- Slide 10: update "@ControllerAdvice" text ---
$s10 = $p.Slides.Item(10)
$shape10 = $s10.Shapes.Item(2)
$shape10.TextFrame.TextRange.Text = "@RestControllerAdvice + @ExceptionHandler"

# --- Slide 5: merge the "都会" + "实现" runs into a single run "都会实现" ---
$s5 = $p.Slides.Item(5)
$shape5 = $s5.Shapes.Item(2)
$para2 = $shape5.TextFrame.TextRange.Paragraphs(2)
# Remove the "实现" run's characters first (cleanly deletes the whole run
# node while offsets still match the original text), then extend the
# preceding "都会" run's text to absorb it.
$para2.Characters(11, 2).Delete()
$para2.Runs(3).Text = "都会实现"
